$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "33.831.45"
$ws.Cells.Item(2, 5).Value2 = "  -0.93%  "
$ws.Cells.Item(3, 4).Value2 = "1.779.34"
$ws.Cells.Item(3, 5).Value2 = "  -1.28%  "
$ws.Cells.Item(4, 5).Value2 = "  +0.11%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = "224.05"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value2 = "  +0.32%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = "0.546"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value2 = "  -1.16%  "
$ws.Cells.Item(7, 5).Value2 = "  +0.12%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = "31.80"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value2 = "  -1.78%  "
$ws.Cells.Item(9, 5).Value2 = "  +0.76%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = "0.0678"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value2 = "  -5.56%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = "0.0935"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value2 = "  +1.02%  "
$ws.Cells.Item(12, 4).Value2 = "2.035.16"
$ws.Cells.Item(12, 5).Value2 = "  -1.30%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value2 = "11.21"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value2 = "  +3.33%  "
$ws.Cells.Item(14, 4).Value2 = "1.764.47"
$ws.Cells.Item(14, 5).Value2 = "  -2.09%  "
$ws.Cells.Item(15, 4).Value2 = "33.859.95"
$ws.Cells.Item(15, 5).Value2 = "  -0.95%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = "0.609"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value2 = "  -3.72%  "
$ws.Cells.Item(17, 5).Value2 = "  -2.29%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = "66.64"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value2 = "  -2.74%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = "238.63"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value2 = "  -3.74%  "
$ws.Cells.Item(20, 4).Value2 = "0.0₃0773"
$ws.Cells.Item(20, 5).Value2 = "  -1.95%  "
$ws.Cells.Item(21, 5).Value2 = "  +0.02%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = "10.56"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value2 = "  -3.40%  "
$ws.Cells.Item(23, 5).Value2 = "  -2.60%  "
$ws.Cells.Item(24, 5).Value2 = "  -2.23%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = "160.79"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value2 = "  +0.50%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = "7.03"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value2 = "  -1.02%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = "16.08"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value2 = "  -3.06%  "
$ws.Cells.Item(28, 5).Value2 = "  -0.91%  "
$ws.Cells.Item(29, 5).Value2 = "  +0.26%  "
$ws.Cells.Item(30, 5).Value2 = "  +0.95%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = "0.0510"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value2 = "  -3.10%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = "3.59"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value2 = "  -3.89%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = "3.51"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value2 = "  -0.20%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = "1.81"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value2 = "  -2.25%  "
$ws.Cells.Item(35, 4).Value2 = "1.388.19"
$ws.Cells.Item(35, 5).Value2 = "  -1.93%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = "0.636"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value2 = "  -2.56%  "
$ws.Cells.Item(37, 5).Value2 = "  -1.81%  "
$ws.Cells.Item(38, 5).Value2 = "  -1.49%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = "2.25"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value2 = "  +4.77%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = "2.38"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value2 = "  +0.73%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = "78.38"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value2 = "  -2.61%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = "0.911"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value2 = "  -3.96%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = "13.50"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value2 = "  +12.30%  "
$ws.Cells.Item(44, 5).Value2 = "  -3.26%  "
$ws.Cells.Item(45, 5).Value2 = "  +11.35%  "
$ws.Cells.Item(46, 5).Value2 = "  +2.59%  "
$ws.Cells.Item(47, 5).Value2 = "  +3.36%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = "5.85"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value2 = "  -1.74%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = "106.49"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value2 = "  -1.09%  "
$ws.Cells.Item(50, 4).Value2 = "1.936.68"
$ws.Cells.Item(50, 5).Value2 = "  -1.34%  "
